$d = $word.ActiveDocument

# 1. Title text change: QTM 385 -> DATASCI 385
$d.Content.Find.Execute("QTM 385 - Experimental Methods", $true, $false, $false, $false, $false,
                         $true, 1, $false, "DATASCI 385 - Experimental Methods", 2)

# 2. Heading text change: "Questions and Answers" -> "Questions"
$d.Content.Find.Execute("Questions and Answers", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Questions", 2)
